$wb = $excel.ActiveWorkbook

# Update the password/value cell on the "loginpage" sheet (B5) to the new value.
$ws1 = $wb.Worksheets.Item("loginpage")
$ws1.Range("B5").Value = "adminn123"

# Make "loginpage" the active sheet / active cell selection (B5), which also
# moves tabSelected off of "managefootertextpage" (previously the active tab).
$ws1.Activate()
$ws1.Range("B5").Select()
